$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40-83 down to 41-84.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new data record.
$ws.Range("A40").Value = 1
$ws.Range("B40").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C40").Value = "Arica y Parinacota"
$ws.Range("D40").Value = 44650
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100102
$ws.Range("H40").Value = "Cítricos"
$ws.Range("I40").Value = 100102005
$ws.Range("J40").Value = "Naranja"
$ws.Range("K40").Value = "Valencia"
$ws.Range("L40").Value = "Segunda"
$ws.Range("M40").Value = 270
$ws.Range("N40").Value = 950
$ws.Range("O40").Value = 1000
$ws.Range("P40").Value = 975
$ws.Range("Q40").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R40").Value = "Región de Coquimbo"
$ws.Range("S40").Value = 975
$ws.Range("T40").Value = 1
